# Auto-generated Excel COM-interop script
# Adds 2024-08-30 violent crime data: updates column K (2024 totals)
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 5232
$ws.Range('K3').Value = 5402
$ws.Range('K4').Value = 1124
$ws.Range('K6').Value = 6015
$ws.Range('K7').Value = 18160

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 159
$ws.Range('K4').Value = 69
$ws.Range('K8').Value = 1226
$ws.Range('K10').Value = 99
$ws.Range('K11').Value = 348
$ws.Range('K14').Value = 95
$ws.Range('K15').Value = 184
$ws.Range('K19').Value = 535
$ws.Range('K20').Value = 420
$ws.Range('K23').Value = 190
$ws.Range('K27').Value = 172
$ws.Range('K29').Value = 970
$ws.Range('K31').Value = 198
$ws.Range('K33').Value = 772
$ws.Range('K34').Value = 103
$ws.Range('K36').Value = 239
$ws.Range('K37').Value = 610
$ws.Range('K42').Value = 668
$ws.Range('K43').Value = 162
$ws.Range('K47').Value = 121
$ws.Range('K48').Value = 226
$ws.Range('K51').Value = 228
$ws.Range('K53').Value = 234
$ws.Range('K54').Value = 356
$ws.Range('K55').Value = 203
$ws.Range('K63').Value = 51
$ws.Range('K65').Value = 414
$ws.Range('K67').Value = 691
$ws.Range('K72').Value = 88
$ws.Range('K76').Value = 251
$ws.Range('K78').Value = 208
$ws.Range('K79').Value = 448
$ws.Range('K83').Value = 403
$ws.Range('K84').Value = 135
$ws.Range('K85').Value = 855
$ws.Range('K86').Value = 122
$ws.Range('K88').Value = 201
$ws.Range('K89').Value = 264
$ws.Range('K90').Value = 163
$ws.Range('K91').Value = 196
$ws.Range('K94').Value = 239
$ws.Range('K95').Value = 311
$ws.Range('K96').Value = 196
$ws.Range('K99').Value = 308
$ws.Range('K101').Value = 18160

# Sheet: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K6').Value = 35
$ws.Range('K7').Value = 95

# Sheet: West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 61
$ws.Range('K3').Value = 37
$ws.Range('K7').Value = 196

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K6').Value = 118
$ws.Range('K7').Value = 348

# Sheet: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 81
$ws.Range('K7').Value = 264

# Sheet: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K3').Value = 287
$ws.Range('K7').Value = 855

# Sheet: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K4').Value = 11
$ws.Range('K7').Value = 234

# Sheet: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 336
$ws.Range('K7').Value = 1226

# Sheet: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 146
$ws.Range('K7').Value = 403

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 211
$ws.Range('K6').Value = 223
$ws.Range('K7').Value = 772

# Sheet: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K3').Value = 112
$ws.Range('K7').Value = 311

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 203
$ws.Range('K7').Value = 610

# Sheet: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 130
$ws.Range('K7').Value = 414

# Sheet: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 77
$ws.Range('K3').Value = 127
$ws.Range('K7').Value = 308

# Sheet: Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K3').Value = 47
$ws.Range('K6').Value = 71
$ws.Range('K7').Value = 198

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K3').Value = 246
$ws.Range('K7').Value = 691

# Sheet: South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 53
$ws.Range('K7').Value = 135

# Sheet: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 190
$ws.Range('K7').Value = 356

# Sheet: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 349
$ws.Range('K6').Value = 269
$ws.Range('K7').Value = 970

# Sheet: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K6').Value = 112
$ws.Range('K7').Value = 226

# Sheet: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K3').Value = 168
$ws.Range('K7').Value = 535

# Sheet: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 54
$ws.Range('K3').Value = 45
$ws.Range('K7').Value = 251

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 177
$ws.Range('K3').Value = 207
$ws.Range('K7').Value = 668

# Sheet: Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K6').Value = 49
$ws.Range('K7').Value = 99

# Sheet: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 61
$ws.Range('K6').Value = 75
$ws.Range('K7').Value = 208

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K2').Value = 62
$ws.Range('K3').Value = 57
$ws.Range('K7').Value = 203

# Sheet: Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K3').Value = 69
$ws.Range('K7').Value = 190

# Sheet: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K6').Value = 47
$ws.Range('K7').Value = 196

# Sheet: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 149
$ws.Range('K3').Value = 145
$ws.Range('K7').Value = 448

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 137
$ws.Range('K6').Value = 122
$ws.Range('K7').Value = 420

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K6').Value = 54
$ws.Range('K7').Value = 239

# Sheet: Garfield Ridge
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K3').Value = 27
$ws.Range('K7').Value = 103

# Sheet: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 103
$ws.Range('K7').Value = 239

# Sheet: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K6').Value = 41
$ws.Range('K7').Value = 121

# Sheet: Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 45
$ws.Range('K7').Value = 184

# Sheet: Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 53
$ws.Range('K3').Value = 36
$ws.Range('K4').Value = 15
$ws.Range('K7').Value = 159

# Sheet: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K2').Value = 52
$ws.Range('K7').Value = 201

# Sheet: Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K2').Value = 46
$ws.Range('K3').Value = 41
$ws.Range('K4').Value = 20
$ws.Range('K7').Value = 172

# Sheet: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K4').Value = 53
$ws.Range('K7').Value = 122

# Sheet: Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K3').Value = 51
$ws.Range('K6').Value = 38
$ws.Range('K7').Value = 163

# Sheet: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 65
$ws.Range('K7').Value = 228

# Sheet: Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K4').Value = 23
$ws.Range('K7').Value = 162

# Sheet: Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K3').Value = 24
$ws.Range('K7').Value = 88

# Sheet: Archer Heights
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('K2').Value = 23
$ws.Range('K7').Value = 69
